# Add function to generate exams and solutions:
# correct the true/false solutions for the "1+1=3" questions (tf rows)
# from "T" to "F" (since 1+1=3 is indeed false), and leave the cursor
# selection on D13, matching the author's final editing position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "F"
$ws.Range("D5").Value = "F"
$ws.Range("D11").Value = "F"
$ws.Range("D12").Value = "F"

[void]$ws.Range("D13").Select()
